$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect to make edits, then restore
# protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A13)
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-12 for illustrative purposes only and are subject to change."
$ws.Rows(13).AutoFit()

# Update Weight (D) and Percent Change (E) figures for rows 2-9
$ws.Range("D2").Value = 0.09058892575405596
$ws.Range("E2").Value = -0.02521773040426378

$ws.Range("D3").Value = 0.1065770701668311
$ws.Range("E3").Value = -0.03083085513139538

$ws.Range("D4").Value = 0.1203061428467596
$ws.Range("E4").Value = -0.0220807795065221

$ws.Range("D5").Value = 0.1414818098834224
$ws.Range("E5").Value = -0.02422797927461151

$ws.Range("D6").Value = 0.137935891547892
$ws.Range("E6").Value = -0.01645464601769908

$ws.Range("D7").Value = 0.1472985754523524
$ws.Range("E7").Value = -0.0254421640026482

$ws.Range("D8").Value = 0.1261869273341888
$ws.Range("E8").Value = -0.03053204353083427

$ws.Range("D9").Value = 0.1296246570144979
$ws.Range("E9").Value = -0.02237403163856533

# Total row: Weight stays at 1, only Percent Change updates
$ws.Range("E10").Value = -0.02442483280713093

# Restore sheet protection
$ws.Protect()
